$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typos in the German (informal "du") instructions text (key = INSTRUCTIONS, column DE)
$ws.Range("C28").Value = "Der Rhythmus wird aus vier, acht oder sechzehn Klängen/Tönen bestehen und du hörst jeweils vier Metronomschläge vor und nach dem Rhythmus.\\Deine Aufgabe ist es, den Rhythmus zu hören und dann auf dasjenige Bild der vier Bilder zu klicken, das mit dem Rhythmus übereinstimmt, den du gerade gehört hast.\\ Lass uns das mal üben."

# Fix typos in the German (formal "Sie") instructions text (key = INSTRUCTIONS, column DE_F)
$ws.Range("D28").Value = "Der Rhythmus wird aus vier, acht oder sechzehn Klängen/Tönen bestehen und Sie hören jeweils vier Metronomschläge vor und nach dem Rhythmus.\\Ihre Aufgabe ist es, den Rhythmus zu hören und dann auf dasjenige Bild der vier Bilder zu klicken, das mit dem Rhythmus übereinstimmt, den Sie gerade gehört haben.\\ Lassen Sie uns das mal üben."

# Restore the selection state recorded in the saved workbook
$ws.Range("D29").Select()
